$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and 1h volume change (column E) values
$ws.Range("D2").Value = "69.652.10"
$ws.Range("D3").Value = "3.504.91"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.73"
$ws.Range("E5").Value = "  -1.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "194.85"
$ws.Range("E6").Value = "  +3.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -5.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.649"
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.40"
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000300"
$ws.Range("E12").Value = "  -2.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.50"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").Value = "4.064.82"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "593.94"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.08"
$ws.Range("D17").Value = "69.827.07"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.72"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("E19").Value = "  +2.20%  "
$ws.Range("D20").Value = "3.511.04"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.988"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.19"
$ws.Range("E22").Value = "  +6.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.27"
$ws.Range("E23").Value = "  +2.74%  "
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "101.55"
$ws.Range("E25").Value = "  -4.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.16"
$ws.Range("E26").Value = "  +4.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.85"
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.51"
$ws.Range("E28").Value = "  -1.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.18"
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.06"
$ws.Range("E30").Value = "  +1.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.28"
$ws.Range("E31").Value = "  +2.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.39"
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.13"
$ws.Range("E34").Value = "  -0.32%  "
$ws.Range("D35").Value = "0.0₃0823"
$ws.Range("E35").Value = "  +6.67%  "
$ws.Range("D36").Value = "3.722.78"
$ws.Range("E36").Value = "  +2.55%  "
$ws.Range("E37").Value = "  -1.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.65"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.391"
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.29"
$ws.Range("E41").Value = "  -1.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "484.43"
$ws.Range("E42").Value = "  -4.45%  "
$ws.Range("E43").Value = "  -2.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0453"
$ws.Range("E44").Value = "  -2.03%  "
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("E46").Value = "  -3.15%  "
$ws.Range("E47").Value = "  -1.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.01"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("E49").Value = "  -4.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000247"
$ws.Range("E50").Value = "  +2.93%  "
$ws.Range("E51").Value = "  +10.36%  "
